$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Achicoria" came in (Vega Modelo de Temuco),
# so it gets inserted at the top of the data block (row 111), pushing the
# existing historical rows (111-128) down to (112-129).
$ws.Rows.Item(111).Insert()

$ws.Cells.Item(111, 1).Value2  = 10
$ws.Cells.Item(111, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(111, 3).Value2  = "La Araucanía"
$ws.Cells.Item(111, 4).Value2  = 45127
$ws.Cells.Item(111, 5).Value2  = 9
$ws.Cells.Item(111, 6).Value2  = 100112010
$ws.Cells.Item(111, 7).Value2  = "Achicoria"
$ws.Cells.Item(111, 8).Value2  = "Sin especificar"
$ws.Cells.Item(111, 9).Value2  = "Primera"
$ws.Cells.Item(111, 10).Value2 = 200
$ws.Cells.Item(111, 11).Value2 = 9000
$ws.Cells.Item(111, 12).Value2 = 9000
$ws.Cells.Item(111, 13).Value2 = 9000
$ws.Cells.Item(111, 14).Value2 = "$/caja 18 unidades"
$ws.Cells.Item(111, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(111, 16).Value2 = 500
$ws.Cells.Item(111, 17).Value2 = 18
$ws.Cells.Item(111, 18).Value2 = "Hortaliza"
